$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 1310.25
$ws.Range("I12").Value = 573.125
$ws.Range("K12").Value = 573.125
$ws.Range("M12").Value = -403.125

$ws.Range("H41").Value = 1481
$ws.Range("I41").Value = 1638.125
$ws.Range("J41").Value = 1271.5
$ws.Range("K41").Value = 1638.125
$ws.Range("L41").Value = 1271.5
$ws.Range("M41").Value = -1198.125
$ws.Range("N41").Value = -2151.5

$ws.Range("H53").Value = 1211
$ws.Range("I53").Value = 868.9286
$ws.Range("K53").Value = 868.9286
$ws.Range("M53").Value = -231.9286

$ws.Range("H76").Value = 8266.666999999999
$ws.Range("J76").Value = 12000
$ws.Range("L76").Value = 12000
$ws.Range("N76").Value = -12630

$ws.Range("H79").Value = 8266.666999999999
$ws.Range("J79").Value = 12000
$ws.Range("L79").Value = 12000
$ws.Range("N79").Value = -14184

$ws.Range("H86").Value = 321438620
$ws.Range("I86").Value = 200014060
$ws.Range("K86").Value = 200014060
$ws.Range("M86").Value = -200012937

$ws.Range("H89").Value = 321438620
$ws.Range("I89").Value = 200014060
$ws.Range("K89").Value = 1000070300
$ws.Range("M89").Value = -1000064684

$ws.Range("H118").Value = 1072.7273
$ws.Range("J118").Value = 1250
$ws.Range("L118").Value = 3750
$ws.Range("N118").Value = -7064

$ws.Range("H126").Value = 125259.664
$ws.Range("J126").Value = 125259.664
$ws.Range("L126").Value = 125259.664
$ws.Range("N126").Value = -135139.664

$ws.Range("H132").Value = 5275.077
$ws.Range("I132").Value = 4416
$ws.Range("K132").Value = 13248
$ws.Range("M132").Value = -10718

$ws.Range("H137").Value = 65135.92
$ws.Range("I137").Value = 104311.375
$ws.Range("J137").Value = 2455.2
$ws.Range("K137").Value = 312934.125
$ws.Range("L137").Value = 7365.599999999999
$ws.Range("M137").Value = -310384.125
$ws.Range("N137").Value = -12465.6

$ws.Range("H138").Value = 4095.192
$ws.Range("J138").Value = 4499.8813
$ws.Range("L138").Value = 13499.6439
$ws.Range("N138").Value = -23779.6439

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13825958
$ws.Range("I32").Value = 13728610
$ws.Range("K32").Value = 13728610
$ws.Range("M32").Value = -13728323

$ws.Range("H61").Value = 3008.4
$ws.Range("I61").Value = 2788.3225
$ws.Range("J61").Value = 3766.4443
$ws.Range("K61").Value = 2788.3225
$ws.Range("L61").Value = 3766.4443
$ws.Range("M61").Value = -2576.3225
$ws.Range("N61").Value = -4190.4443

$ws.Range("H63").Value = 6000
$ws.Range("J63").Value = 7400
$ws.Range("L63").Value = 7400
$ws.Range("N63").Value = -8772

$ws.Range("H66").Value = 6000
$ws.Range("J66").Value = 7400
$ws.Range("L66").Value = 37000
$ws.Range("N66").Value = -43864

$ws.Range("H74").Value = 2163.56
$ws.Range("I74").Value = 2154.5
$ws.Range("J74").Value = 2199.8
$ws.Range("K74").Value = 2154.5
$ws.Range("L74").Value = 2199.8
$ws.Range("M74").Value = -1280.5
$ws.Range("N74").Value = -3947.8

$ws.Range("H77").Value = 2163.56
$ws.Range("I77").Value = 2154.5
$ws.Range("J77").Value = 2199.8
$ws.Range("K77").Value = 10772.5
$ws.Range("L77").Value = 10999
$ws.Range("M77").Value = -6404.5
$ws.Range("N77").Value = -19735

$ws.Range("H109").Value = 35000
$ws.Range("J109").Value = 35000
$ws.Range("L109").Value = 35000
$ws.Range("N109").Value = -37774

$ws.Range("H114").Value = 107824.5
$ws.Range("J114").Value = 107824.5
$ws.Range("L114").Value = 107824.5
$ws.Range("N114").Value = -116502.5

$ws.Range("H136").Value = 3008.4
$ws.Range("I136").Value = 2788.3225
$ws.Range("J136").Value = 3766.4443
$ws.Range("K136").Value = 8364.967500000001
$ws.Range("L136").Value = 11299.3329
$ws.Range("M136").Value = -5814.967500000001
$ws.Range("N136").Value = -16399.3329

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()

$ws.Range("H134").Value = 1077174.5
$ws.Range("I134").Value = 1171298.9
$ws.Range("K134").Value = 3513896.7
$ws.Range("M134").Value = -3511361.7

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 9968.333000000001
$ws.Range("J22").Value = 22175.6
$ws.Range("L22").Value = 22175.6
$ws.Range("N22").Value = -22875.6

$ws.Range("H28").Value = 58501.5
$ws.Range("J28").Value = 58501.5
$ws.Range("L28").Value = 58501.5
$ws.Range("N28").Value = -58991.5

$ws.Range("H132").Value = 3005.7817
$ws.Range("I132").Value = 2884.0222
$ws.Range("K132").Value = 8652.0666
$ws.Range("M132").Value = -6122.0666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 3010.0881
$ws.Range("J68").Value = 3390.28
$ws.Range("L68").Value = 10170.84
$ws.Range("N68").Value = -11792.84

$ws.Range("H71").Value = 3010.0881
$ws.Range("J71").Value = 3390.28
$ws.Range("L71").Value = 30512.52
$ws.Range("N71").Value = -38624.52

$ws.Range("H107").Value = 1199.0217
$ws.Range("I107").Value = 986.2143
$ws.Range("K107").Value = 2958.6429
$ws.Range("M107").Value = -1038.6429

$ws.Range("H114").Value = 977.375
$ws.Range("J114").Value = 415.6
$ws.Range("L114").Value = 1246.8
$ws.Range("N114").Value = -7754.8

$ws.Range("H117").Value = 169398.83
$ws.Range("J117").Value = 169398.83
$ws.Range("L117").Value = 508196.49
$ws.Range("N117").Value = -515080.49

$ws.Range("H121").Value = 10050
$ws.Range("J121").Value = 20000
$ws.Range("L121").Value = 60000
$ws.Range("N121").Value = -62620

$ws.Range("H131").Value = 1658.4166
$ws.Range("J131").Value = 1753.7059
$ws.Range("L131").Value = 5261.1177
$ws.Range("N131").Value = -15341.1177

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2800.625
$ws.Range("I22").Value = 2557.8572
$ws.Range("K22").Value = 2557.8572
$ws.Range("M22").Value = -2262.8572

$ws.Range("H27").Value = 2800.625
$ws.Range("I27").Value = 2557.8572
$ws.Range("K27").Value = 2557.8572
$ws.Range("M27").Value = -2450.8572

$ws.Range("H40").Value = 83335064
$ws.Range("I40").Value = 111112750
$ws.Range("K40").Value = 111112750
$ws.Range("M40").Value = -111112614

$ws.Range("H46").Value = 2980.25
$ws.Range("I46").Value = 2055.5715
$ws.Range("K46").Value = 2055.5715
$ws.Range("M46").Value = -1867.5715

$ws.Range("H122").Value = 17636.273
$ws.Range("J122").Value = 8749.75
$ws.Range("L122").Value = 26249.25
$ws.Range("N122").Value = -31149.25

$ws.Range("H136").Value = 4935.8486
$ws.Range("I136").Value = 4883.625
$ws.Range("K136").Value = 14650.875
$ws.Range("M136").Value = -12100.875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 470393.62
$ws.Range("I2").Value = 17433
$ws.Range("K2").Value = 17433
$ws.Range("M2").Value = -17321

$ws.Range("H4").Value = 885960.9399999999
$ws.Range("J4").Value = 1252002.5
$ws.Range("L4").Value = 1252002.5
$ws.Range("N4").Value = -1252228.5

$ws.Range("H31").Value = 39998.5
$ws.Range("J31").Value = 39998.5
$ws.Range("L31").Value = 39998.5
$ws.Range("N31").Value = -40694.5

$ws.Range("H51").Value = 59998
$ws.Range("I51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("M51").ClearContents()

$ws.Range("H130").Value = 53214.5
$ws.Range("J130").Value = 53214.5
$ws.Range("L130").Value = 53214.5
$ws.Range("N130").Value = -63254.5

$ws.Range("H136").Value = 21386.818
$ws.Range("I136").Value = 3122.9143
$ws.Range("J136").Value = 53348.65
$ws.Range("K136").Value = 9368.742899999999
$ws.Range("L136").Value = 160045.95
$ws.Range("M136").Value = -6818.742899999999
$ws.Range("N136").Value = -165145.95
